# "Arregladas fechas y agregado un .exe"
#
# The two date cells on Hoja1 (D3 = fecha_ingreso, H3 = fecha_cargo) were
# re-entered as literal text (dd/mm/yyyy) instead of real Excel date
# serials, so the whole D and H data columns were switched to a Text
# number format to stop Excel from reinterpreting them back into dates.
# Cell E6 picked up an underline, and the saved selection moved to H3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Column D (fecha_ingreso) and column H (fecha_cargo), rows 3-12: switch to
# Text format ("@") so the dates stick as plain strings.
$ws.Range("D3:D12").NumberFormat = "@"
$ws.Range("H3:H12").NumberFormat = "@"

# Re-type the two existing dates as literal text.
$ws.Range("D3").Value = "09/08/1997"
$ws.Range("H3").Value = "01/08/2021"

# Underline cell E6.
$ws.Range("E6").Font.Underline = $true

# Leave the saved cursor/selection on H3, matching the workbook view state.
$ws.Range("H3").Select()

$wb.Save()
